$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue "D2" "98.748.77"
Set-TextValue "E2" "  +2.33%  "
Set-TextValue "D3" "3.318.45"
Set-TextValue "E3" "  +1.35%  "
Set-TextValue "E4" "  -0.03%  "
Set-TextValue "D5" "256.34"
Set-TextValue "E5" "  +2.87%  "
Set-TextValue "D6" "624.89"
Set-TextValue "E6" "  +1.32%  "
Set-TextValue "E7" "  +29.77%  "
Set-TextValue "D8" "0.405"
Set-TextValue "E8" "  +6.51%  "
Set-TextValue "E9" "  -0.09%  "
Set-TextValue "D10" "0.924"
Set-TextValue "E10" "  +18.73%  "
Set-TextValue "D11" "3.315.59"
Set-TextValue "E11" "  +1.46%  "
Set-TextValue "D12" "0.200"
Set-TextValue "E12" "  +0.91%  "
Set-TextValue "D13" "39.35"
Set-TextValue "E13" "  +12.60%  "
Set-TextValue "D14" "98.404.92"
Set-TextValue "E14" "  +2.33%  "
Set-TextValue "D15" "0.0000250"
Set-TextValue "E15" "  +2.87%  "
Set-TextValue "D16" "3.944.22"
Set-TextValue "E16" "  +1.79%  "
Set-TextValue "D17" "5.50"
Set-TextValue "E17" "  +0.45%  "
Set-TextValue "D18" "3.317.92"
Set-TextValue "E18" "  +1.46%  "
Set-TextValue "D19" "3.50"
Set-TextValue "E19" "  -1.67%  "
Set-TextValue "D20" "15.36"
Set-TextValue "E20" "  +3.70%  "
Set-TextValue "D21" "6.33"
Set-TextValue "E21" "  +9.80%  "
Set-TextValue "D22" "484.78"
Set-TextValue "E22" "  +1.79%  "
Set-TextValue "D23" "9.46"
Set-TextValue "E23" "  +3.47%  "
Set-TextValue "E24" "  +1.25%  "
Set-TextValue "D25" "5.63"
Set-TextValue "E25" "  +1.08%  "
Set-TextValue "D26" "88.54"
Set-TextValue "E26" "  +1.32%  "
Set-TextValue "D27" "12.01"
Set-TextValue "E27" "  +0.52%  "
Set-TextValue "D28" "3.490.16"
Set-TextValue "E28" "  +1.16%  "
Set-TextValue "D29" "0.292"
Set-TextValue "E29" "  +23.52%  "
Set-TextValue "E30" "  -0.20%  "
Set-TextValue "D31" "0.188"
Set-TextValue "E31" "  +4.11%  "
Set-TextValue "D32" "0.130"
Set-TextValue "E32" "  +8.39%  "
Set-TextValue "D33" "10.26"
Set-TextValue "E33" "  +12.26%  "
Set-TextValue "E34" "  +0.13%  "
Set-TextValue "D35" "27.97"
Set-TextValue "E35" "  +3.35%  "
Set-TextValue "D36" "7.25"
Set-TextValue "E36" "  -0.91%  "
Set-TextValue "E37" "  -0.45%  "
Set-TextValue "E38" "  +1.69%  "
Set-TextValue "E39" "  +4.76%  "
Set-TextValue "E40" "  +0.51%  "
Set-TextValue "D41" "492.03"
Set-TextValue "E41" "  +0.37%  "
Set-TextValue "E42" "  -1.27%  "
Set-TextValue "D43" "3.62"
Set-TextValue "E43" "  +5.50%  "
Set-TextValue "D44" "0.797"
Set-TextValue "E44" "  +1.76%  "
Set-TextValue "E45" "  -0.02%  "
Set-TextValue "D46" "3.15"
Set-TextValue "E46" "  -2.34%  "
Set-TextValue "D47" "158.74"
Set-TextValue "E47" "  -1.44%  "
Set-TextValue "D48" "7.45"
Set-TextValue "E48" "  +18.43%  "
Set-TextValue "D49" "1.94"
Set-TextValue "E49" "  +2.98%  "
Set-TextValue "D50" "0.845"
Set-TextValue "E50" "  +8.67%  "
Set-TextValue "D51" "4.70"
Set-TextValue "E51" "  +5.34%  "
